# Update the dSF column (column F) values on the data sheet
# to reflect the repulled data / recalculated mean values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value  = -2
$ws.Range("F6").Value  = -3
$ws.Range("F9").Value  = -4
$ws.Range("F13").Value = 3
$ws.Range("F15").Value = -3
$ws.Range("F18").Value = -2
$ws.Range("F26").Value = -1
$ws.Range("F31").Value = -2
$ws.Range("F33").Value = -6
